$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.6077046446798136
$ws.Range("D2").Value = 0.03079724402117989
$ws.Range("E2").Value = 0.224429308261251
$ws.Range("F2").Value = 0.7248687392599535
$ws.Range("G2").Value = 0.00242972007515663
$ws.Range("I2").Value = 0.9809152544316149
$ws.Range("K2").Value = 0.4690150215434414
$ws.Range("L2").Value = 0.1719004026205795
$ws.Range("O2").Value = 2.497026713160096
$ws.Range("B3").Value = 0.5796414754831005
$ws.Range("D3").Value = 0.02879354639640752
$ws.Range("E3").Value = 0.2265548982605265
$ws.Range("F3").Value = 0.7238108711795519
$ws.Range("G3").Value = 0.002432198762770559
$ws.Range("I3").Value = 0.9927891755959593
$ws.Range("K3").Value = 0.4104088994488677
$ws.Range("L3").Value = 0.1601610097961412
$ws.Range("O3").Value = 2.507101683807377
$ws.Range("B4").Value = 0.5626095196223559
$ws.Range("D4").Value = 0.02755315446407991
$ws.Range("E4").Value = 0.2279470801253858
$ws.Range("F4").Value = 0.723645887686537
$ws.Range("G4").Value = 0.002433802736840874
$ws.Range("I4").Value = 1.00057301700701
$ws.Range("K4").Value = 0.3742607026897247
$ws.Range("L4").Value = 0.1530203475481215
$ws.Range("O4").Value = 2.514886927362753
$ws.Range("B5").Value = 0.5557194436901511
$ws.Range("D5").Value = 0.02704516848562832
$ws.Range("E5").Value = 0.2285363285517263
$ws.Range("F5").Value = 0.723700466314952
$ws.Range("G5").Value = 0.002434477062480316
$ws.Range("I5").Value = 1.003869013700236
$ws.Range("K5").Value = 0.3594896551382476
$ws.Range("L5").Value = 0.1501275280282499
$ws.Range("O5").Value = 2.518461400828329
$ws.Range("B6").Value = 0.5545784234709572
$ws.Range("D6").Value = 0.026960666623566
$ws.Range("E6").Value = 0.2286354978052074
$ws.Range("F6").Value = 0.7237168851329017
$ws.Range("G6").Value = 0.002434590285278813
$ws.Range("I6").Value = 1.004423801757374
$ws.Range("K6").Value = 0.3570345204384751
$ws.Range("L6").Value = 0.1496482112839317
$ws.Range("O6").Value = 2.519079208723497
$ws.Range("B7").Value = 0.5625163920921636
$ws.Range("D7").Value = 0.02754631373974803
$ws.Range("E7").Value = 0.2279549381287982
$ws.Range("F7").Value = 0.7236461305918525
$ws.Range("G7").Value = 0.002433811747382583
$ws.Range("I7").Value = 1.000616965830133
$ws.Range("K7").Value = 0.3740616575116746
$ws.Range("L7").Value = 0.1529812647481634
$ws.Range("O7").Value = 2.514933506859919
$ws.Range("B8").Value = 0.5979874980024249
$ws.Range("D8").Value = 0.03010848370041685
$ws.Range("E8").Value = 0.2251441608267335
$ws.Range("F8").Value = 0.7244033949001434
$ws.Range("G8").Value = 0.002430557733685446
$ws.Range("I8").Value = 0.9849070170767149
$ws.Range("K8").Value = 0.4488421042360642
$ws.Range("L8").Value = 0.1678387431039283
$ws.Range("O8").Value = 2.500168529963247
$ws.Range("B9").Value = 0.6691041461781992
$ws.Range("D9").Value = 0.03505168422596938
$ws.Range("E9").Value = 0.2203216144610414
$ws.Range("F9").Value = 0.7297358508858522
$ws.Range("G9").Value = 0.002424824875068473
$ws.Range("I9").Value = 0.9580130718113153
$ws.Range("K9").Value = 0.5941572759948599
$ws.Range("L9").Value = 0.197505238065574
$ws.Range("O9").Value = 2.483914633863122
$ws.Range("B10").Value = 0.7222813739776086
$ws.Range("D10").Value = 0.03863298298050211
$ws.Range("E10").Value = 0.2171967730951474
$ws.Range("F10").Value = 0.7360047355786676
$ws.Range("G10").Value = 0.002421004242896591
$ws.Range("I10").Value = 0.9406391104756686
$ws.Range("K10").Value = 0.7000807689378519
$ws.Range("L10").Value = 0.2196225295919447
$ws.Range("O10").Value = 2.479734630999531
$ws.Range("B11").Value = 0.7466702662102875
$ws.Range("D11").Value = 0.04025106503587494
$ws.Range("E11").Value = 0.2158656061615449
$ws.Range("F11").Value = 0.7393684441233148
$ws.Range("G11").Value = 0.002419350277192517
$ws.Range("I11").Value = 0.9332530027675858
$ws.Range("K11").Value = 0.7480802438620344
$ws.Range("L11").Value = 0.2297536602258674
$ws.Range("O11").Value = 2.479522802189535
$ws.Range("B12").Value = 0.7559337375154485
$ws.Range("D12").Value = 0.04086217595855857
$ws.Range("E12").Value = 0.2153744871477894
$ws.Range("F12").Value = 0.7407158809475476
$ws.Range("G12").Value = 0.002418735988704661
$ws.Range("I12").Value = 0.9305304840402293
$ws.Range("K12").Value = 0.766229006172523
$ws.Range("L12").Value = 0.2336000165675785
$ws.Range("O12").Value = 2.479685853896541
$ws.Range("B13").Value = 0.7539374505172134
$ws.Range("D13").Value = 0.04073063488630879
$ws.Range("E13").Value = 0.2154796821996587
$ws.Range("F13").Value = 0.7404224088321172
$ws.Range("G13").Value = 0.002418867752461246
$ws.Range("I13").Value = 0.9311135157979713
$ws.Range("K13").Value = 0.7623215871700779
$ws.Range("L13").Value = 0.232771195907219
$ws.Range("O13").Value = 2.47963991380584
$ws.Range("B14").Value = 0.7474318212359208
$ws.Range("D14").Value = 0.0403013741705962
$ws.Range("E14").Value = 0.215824941812981
$ws.Range("F14").Value = 0.739477821733729
$ws.Range("G14").Value = 0.002419299498410308
$ws.Range("I14").Value = 0.9330275274008386
$ws.Range("K14").Value = 0.7495739116913001
$ws.Range("L14").Value = 0.2300699042654912
$ws.Range("O14").Value = 2.479531339456457
$ws.Range("B15").Value = 0.7434505549437063
$ws.Range("D15").Value = 0.04003822751494823
$ws.Range("E15").Value = 0.2160381109636127
$ws.Range("F15").Value = 0.7389088308872189
$ws.Range("G15").Value = 0.002419565520702217
$ws.Range("I15").Value = 0.9342096102190816
$ws.Range("K15").Value = 0.7417619717631112
$ws.Range("L15").Value = 0.2284165727883902
$ws.Range("O15").Value = 2.479496523686038
$ws.Range("B16").Value = 0.720691442533564
$ws.Range("D16").Value = 0.03852701221816801
$ws.Range("E16").Value = 0.2172855836446175
$ws.Range("F16").Value = 0.7357952159377277
$ws.Range("G16").Value = 0.002421114019991691
$ws.Range("I16").Value = 0.9411322246412261
$ws.Range("K16").Value = 0.6969400774898986
$ws.Range("L16").Value = 0.2189618311710575
$ws.Range("O16").Value = 2.479782500257272
$ws.Range("B17").Value = 0.7067798352942987
$ws.Range("D17").Value = 0.03759707427100523
$ws.Range("E17").Value = 0.2180739874645781
$ws.Range("F17").Value = 0.734016278534618
$ws.Range("G17").Value = 0.002422085464290635
$ws.Range("I17").Value = 0.9455115694426439
$ws.Range("K17").Value = 0.6693951049745124
$ws.Range("L17").Value = 0.2131794553083211
$ws.Range("O17").Value = 2.480390913859026
$ws.Range("B18").Value = 0.6987969472156124
$ws.Range("D18").Value = 0.0370611585059919
$ws.Range("E18").Value = 0.2185359611591737
$ws.Range("F18").Value = 0.7330412658994661
$ws.Range("G18").Value = 0.002422652128483693
$ws.Range("I18").Value = 0.9480791500081622
$ws.Range("K18").Value = 0.653534544326817
$ws.Range("L18").Value = 0.2098601693055002
$ws.Range("O18").Value = 2.480899880178271
$ws.Range("B19").Value = 0.6960973082686621
$ws.Range("D19").Value = 0.03687952892041579
$ws.Range("E19").Value = 0.2186938389816913
$ws.Range("F19").Value = 0.732719417543926
$ws.Range("G19").Value = 0.002422845352507499
$ws.Range("I19").Value = 0.9489568512774227
$ws.Range("K19").Value = 0.6481614664405697
$ws.Range("L19").Value = 0.2087374507822233
$ws.Range("O19").Value = 2.48109951077879
$ws.Range("B20").Value = 0.708258818477475
$ws.Range("D20").Value = 0.03769617561173533
$ws.Range("E20").Value = 0.2179891805169953
$ws.Range("F20").Value = 0.7342006620835377
$ws.Range("G20").Value = 0.00242198123357961
$ws.Range("I20").Value = 0.9450403404429366
$ws.Range("K20").Value = 0.6723291238643583
$ws.Range("L20").Value = 0.2137943186433517
$ws.Range("O20").Value = 2.480309686778867
$ws.Range("B21").Value = 0.749341929004629
$ws.Range("D21").Value = 0.04042750269758244
$ws.Range("E21").Value = 0.2157231790325653
$ws.Range("F21").Value = 0.7397532699646803
$ws.Range("G21").Value = 0.002419172358122305
$ws.Range("I21").Value = 0.9324633152046964
$ws.Range("K21").Value = 0.7533189691220628
$ws.Range("L21").Value = 0.2308630709873256
$ws.Range("O21").Value = 2.479556625864973
$ws.Range("B22").Value = 0.7763546252447782
$ws.Range("D22").Value = 0.04220311685978828
$ws.Range("E22").Value = 0.2143177742601754
$ws.Range("F22").Value = 0.7438116690673695
$ws.Range("G22").Value = 0.002417406703488663
$ws.Range("I22").Value = 0.9246774134393014
$ws.Range("K22").Value = 0.806089242626058
$ws.Range("L22").Value = 0.2420762313106763
$ws.Range("O22").Value = 2.480482535226059
$ws.Range("B23").Value = 0.761922765459758
$ws.Range("D23").Value = 0.04125631418306597
$ws.Range("E23").Value = 0.2150609601035143
$ws.Range("F23").Value = 0.7416063130723956
$ws.Range("G23").Value = 0.002418342669901818
$ws.Range("I23").Value = 0.9287931798602571
$ws.Range("K23").Value = 0.7779398074369226
$ws.Range("L23").Value = 0.2360863158178574
$ws.Range("O23").Value = 2.479858509640565
$ws.Range("B24").Value = 0.7075901231632713
$ws.Range("D24").Value = 0.03765137588689527
$ws.Range("E24").Value = 0.2180274946063934
$ws.Range("F24").Value = 0.7341171536336333
$ws.Range("G24").Value = 0.002422028330952721
$ws.Range("I24").Value = 0.9452532278374441
$ws.Range("K24").Value = 0.671002730450482
$ws.Range("L24").Value = 0.2135163231009756
$ws.Range("O24").Value = 2.480345913721493
$ws.Range("B25").Value = 0.6497006961967031
$ws.Range("D25").Value = 0.03372321134861522
$ws.Range("E25").Value = 0.2215526433365254
$ws.Range("F25").Value = 0.7278807974845449
$ws.Range("G25").Value = 0.002426306769294988
$ws.Range("I25").Value = 0.9648697203492382
$ws.Range("K25").Value = 0.554991004872619
$ws.Range("L25").Value = 0.1678387431039283
$ws.Range("O25").Value = 2.486950090088385
